# Oprava termínu školení v leafletu.
# The "Kdy:" box on slide 1 currently reads "Kdy:  3-4. října 2014".
# Split the date run into three runs so the month can read "listopadu"
# instead of "října", keeping the existing run formatting (same rPr).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the shape that holds the "Kdy:" / date text (works even if shape
# ordering/ids ever shift) instead of hard-coding an index.
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -like "Kdy:*") {
            $targetShape = $shp
            break
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$datePara = $tr.Paragraphs(1, 1)

# Run 1 = "Kdy:" (bold), Run 2 = "  3-4. října 2014" (regular).
$dateRun = $datePara.Runs(2, 1)

# Shrink the existing run down to "  3-4. " and append two new runs that
# inherit its formatting, producing: "  3-4. " + "listopadu " + "2014".
$dateRun.Text = "  3-4. "
$monthRun = $dateRun.InsertAfter("listopadu ")
$yearRun = $monthRun.InsertAfter("2014")
